# evaluation/Versuchsplan.xlsx — "update bias calc for no edges"
#
# D455 ("Bias") sheet gets 5 new per-edge diagnostic columns (NaN-Ratio,
# Edge left/down/right/up) inserted between "Precision" and "ADR", with
# sample data filled in for the first two rows. All sheets get their zoom
# dropped to 90%, and the active sheet/tab moves from "Astra Stereo" to
# "D455".

$wb = $excel.ActiveWorkbook

$ws435   = $wb.Worksheets.Item(1)   # D435
$ws455   = $wb.Worksheets.Item(2)   # D455
$wsZed   = $wb.Worksheets.Item(3)   # ZED2
$wsOakD  = $wb.Worksheets.Item(4)   # OAK-D
$wsOakDP = $wb.Worksheets.Item(5)   # OAK-D Pro
$wsAstra = $wb.Worksheets.Item(6)   # Astra Stereo

# --- D455: shift "ADR"/"Sphere" right and insert the 5 new headers -------
# Old layout: ... I=Bias J=Precision K=Edge Precision L=ADR   M=Sphere
# New layout: ... I=Bias J=Precision K=NaN-Ratio L=Edge left M=Edge down
#             N=Edge right O=Edge up P=ADR Q=Sphere
$ws455.Range("Q1").Value = $ws455.Range("M1").Value()
$ws455.Range("P1").Value = $ws455.Range("L1").Value()
$ws455.Range("K1").Value = "NaN-Ratio"
$ws455.Range("L1").Value = "Edge left"
$ws455.Range("M1").Value = "Edge down"
$ws455.Range("N1").Value = "Edge right"
$ws455.Range("O1").Value = "Edge up"

# New per-row diagnostic data (rows 2 and 3)
$ws455.Range("I2").Value = 0.00452
$ws455.Range("J2").Value = 0.00528
$ws455.Range("K2").Value = 0.00048
$ws455.Range("L2").Value = 1.75261
$ws455.Range("M2").Value = 0.83784
$ws455.Range("N2").Value = 0.78911
$ws455.Range("O2").Value = 1.62116

$ws455.Range("I3").Value = 0.01079
$ws455.Range("J3").Value = 0.021
$ws455.Range("K3").Value = 0.00237
$ws455.Range("L3").Value = 3.20879
$ws455.Range("M3").Value = 4.44525
$ws455.Range("N3").Value = 2.79416
$ws455.Range("O3").Value = 1.75754

# Widen the two new "Edge left"/"Edge down" label columns (K:L) so the
# headers aren't clipped, matching the rest of the sheet's wide columns.
$ws455.Columns.Item(11).ColumnWidth = 19.25
$ws455.Columns.Item(12).ColumnWidth = 19.25

# Clean out the legacy trailing blank rows below the table (rows
# 1048563:1048576) that used to pad the sheet out.
$ws455.Range("A1048563:A1048576").EntireRow.Delete()

# Move the cursor/selection to the bottom-right of the new table.
$ws455.Range("L18").Select()

# --- Zoom: every sheet drops from 100% to 90% -----------------------------
$ws435.Activate()
$excel.ActiveWindow.Zoom = 90
$ws455.Activate()
$excel.ActiveWindow.Zoom = 90
$wsZed.Activate()
$excel.ActiveWindow.Zoom = 90
$wsOakD.Activate()
$excel.ActiveWindow.Zoom = 90
$wsOakDP.Activate()
$excel.ActiveWindow.Zoom = 90
$wsAstra.Activate()
$excel.ActiveWindow.Zoom = 90

# --- Active tab moves from "Astra Stereo" (index 5) to "D455" (index 1) --
$ws455.Activate()
